$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2013
$ws.Range("I28").Value = 592.3333
$ws.Range("K28").Value = 592.3333
$ws.Range("M28").Value = -107.3333
$ws.Range("H40").Value = 3852.923
$ws.Range("I40").Value = 2939.2
$ws.Range("J40").Value = 4424
$ws.Range("K40").Value = 2939.2
$ws.Range("L40").Value = 4424
$ws.Range("M40").Value = -2764.2
$ws.Range("N40").Value = -4774
$ws.Range("H64").Value = 5046.5
$ws.Range("I64").Value = 4744
$ws.Range("J64").Value = 5147.3335
$ws.Range("K64").Value = 4744
$ws.Range("L64").Value = 5147.3335
$ws.Range("M64").Value = -4496
$ws.Range("N64").Value = -5643.3335
$ws.Range("H67").Value = 5046.5
$ws.Range("I67").Value = 4744
$ws.Range("J67").Value = 5147.3335
$ws.Range("K67").Value = 4744
$ws.Range("L67").Value = 5147.3335
$ws.Range("M67").Value = -3886
$ws.Range("N67").Value = -6863.3335
$ws.Range("H74").Value = 8171.1665
$ws.Range("I74").Value = 6594.4287
$ws.Range("J74").Value = 10378.6
$ws.Range("K74").Value = 6594.4287
$ws.Range("L74").Value = 10378.6
$ws.Range("M74").Value = -5658.4287
$ws.Range("N74").Value = -12250.6
$ws.Range("H77").Value = 8171.1665
$ws.Range("I77").Value = 6594.4287
$ws.Range("J77").Value = 10378.6
$ws.Range("K77").Value = 32972.14350000001
$ws.Range("L77").Value = 51893
$ws.Range("M77").Value = -28292.14350000001
$ws.Range("N77").Value = -61253
$ws.Range("H107").Value = 499
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 499
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 499
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -4339
$ws.Range("H112").Value = 3572
$ws.Range("J112").Value = 4054.7778
$ws.Range("L112").Value = 12164.3334
$ws.Range("N112").Value = -14380.3334
$ws.Range("H138").Value = 4311.5264
$ws.Range("J138").Value = 3892.25
$ws.Range("L138").Value = 11676.75
$ws.Range("N138").Value = -21956.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14422.3955
$ws.Range("I32").Value = 14845.411
$ws.Range("J32").Value = 10562.375
$ws.Range("K32").Value = 14845.411
$ws.Range("L32").Value = 10562.375
$ws.Range("M32").Value = -14558.411
$ws.Range("N32").Value = -11136.375
$ws.Range("H45").Value = 4019.5334
$ws.Range("I45").Value = 3581.3635
$ws.Range("J45").Value = 5224.5
$ws.Range("K45").Value = 3581.3635
$ws.Range("L45").Value = 5224.5
$ws.Range("M45").Value = -3204.3635
$ws.Range("N45").Value = -5978.5
$ws.Range("H61").Value = 6953.931
$ws.Range("I61").Value = 6827.2856
$ws.Range("K61").Value = 6827.2856
$ws.Range("M61").Value = -6615.2856
$ws.Range("H63").Value = 5049.8887
$ws.Range("I63").Value = 4431.375
$ws.Range("K63").Value = 4431.375
$ws.Range("M63").Value = -3745.375
$ws.Range("H66").Value = 5049.8887
$ws.Range("I66").Value = 4431.375
$ws.Range("K66").Value = 22156.875
$ws.Range("M66").Value = -18724.875
$ws.Range("H74").Value = 2728.6572
$ws.Range("I74").Value = 2505.6453
$ws.Range("J74").Value = 4457
$ws.Range("K74").Value = 2505.6453
$ws.Range("L74").Value = 4457
$ws.Range("M74").Value = -1631.6453
$ws.Range("N74").Value = -6205
$ws.Range("H77").Value = 2728.6572
$ws.Range("I77").Value = 2505.6453
$ws.Range("J77").Value = 4457
$ws.Range("K77").Value = 12528.2265
$ws.Range("L77").Value = 22285
$ws.Range("M77").Value = -8160.226500000001
$ws.Range("N77").Value = -31021
$ws.Range("H132").Value = 56441.05
$ws.Range("I132").Value = 62493.117
$ws.Range("K132").Value = 187479.351
$ws.Range("M132").Value = -184949.351
$ws.Range("H136").Value = 6953.931
$ws.Range("I136").Value = 6827.2856
$ws.Range("K136").Value = 20481.8568
$ws.Range("M136").Value = -17931.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2776.7568
$ws.Range("I86").Value = 1158
$ws.Range("J86").Value = 3762.087
$ws.Range("K86").Value = 1158
$ws.Range("L86").Value = 3762.087
$ws.Range("M86").Value = -35
$ws.Range("N86").Value = -6008.087
$ws.Range("H89").Value = 2776.7568
$ws.Range("I89").Value = 1158
$ws.Range("J89").Value = 3762.087
$ws.Range("K89").Value = 5790
$ws.Range("L89").Value = 18810.435
$ws.Range("M89").Value = -174
$ws.Range("N89").Value = -30042.435
$ws.Range("H94").Value = 2354.75
$ws.Range("I94").Value = 2042.3334
$ws.Range("J94").Value = 5166.5
$ws.Range("K94").Value = 2042.3334
$ws.Range("L94").Value = 5166.5
$ws.Range("M94").Value = -1591.3334
$ws.Range("N94").Value = -6068.5
$ws.Range("H103").Value = 14625.5
$ws.Range("J103").Value = 14625.5
$ws.Range("L103").Value = 14625.5
$ws.Range("N103").Value = -16969.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 548.4074000000001
$ws.Range("I22").Value = 468.5
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 468.5
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -118.5
$ws.Range("N22").Value = -1600
$ws.Range("H105").Value = 1037.1111
$ws.Range("I105").Value = 1037.1111
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1037.1111
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 709.8888999999999
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 348.75
$ws.Range("I81").Value = 348.75
$ws.Range("K81").Value = 1046.25
$ws.Range("M81").Value = 76.75
$ws.Range("H84").Value = 348.75
$ws.Range("I84").Value = 348.75
$ws.Range("K84").Value = 3138.75
$ws.Range("M84").Value = 2477.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 38368.484
$ws.Range("I132").Value = 39524.5
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 118573.5
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -116043.5
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 15420.444
$ws.Range("I46").Value = 25429.8
$ws.Range("J46").Value = 2908.75
$ws.Range("K46").Value = 25429.8
$ws.Range("L46").Value = 2908.75
$ws.Range("M46").Value = -25241.8
$ws.Range("N46").Value = -3284.75
$ws.Range("H101").Value = 32000
$ws.Range("J101").Value = 32000
$ws.Range("L101").Value = 32000
$ws.Range("N101").Value = -38490
$ws.Range("H136").Value = 3666.818
$ws.Range("I136").Value = 2285.2144
$ws.Range("J136").Value = 6084.625
$ws.Range("K136").Value = 6855.6432
$ws.Range("L136").Value = 18253.875
$ws.Range("M136").Value = -4305.6432
$ws.Range("N136").Value = -23353.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = ""
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = ""
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""
$ws.Range("H136").Value = 6229.25
$ws.Range("I136").Value = 7668
$ws.Range("J136").Value = 5749.6665
$ws.Range("K136").Value = 23004
$ws.Range("L136").Value = 17248.9995
$ws.Range("M136").Value = -20454
$ws.Range("N136").Value = -22348.9995

Write-Output "Louisoix_Profits values updated"

